$d = $word.ActiveDocument

# Locate the "|7| - Lich" paragraph (end of the "Parts:" legend list). It is
# immediately followed by a blank separator paragraph and then "EXAMPLE:".
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq "|7| - Lich") {
        $anchorIndex = $i
        break
    }
}

# Insert a new blank paragraph right after "|7| - Lich" (formatting is
# cloned from the anchor paragraph, matching the rest of the list), then
# append four more paragraphs carrying the new legend lines, each cloning
# the same paragraph/run formatting as its predecessor.
$insertAt = $anchorIndex + 1
$d.Paragraphs($anchorIndex).Range.InsertParagraphAfter()

$texts = @(
    "|8| - Necromancer fight line",
    "|9| - Regular bandit line",
    "|10| - Not belonging to this world",
    "|11| - Bandit leader line"
)

foreach ($t in $texts) {
    $d.Paragraphs($insertAt).Range.InsertParagraphAfter()
    $insertAt = $insertAt + 1
    $d.Paragraphs($insertAt).Range.InsertBefore($t)
}
